$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "X2"  = 0.05591268251082138
    "Y2"  = 0.4494154161992178
    "AB2" = 0.05618208824592539
    "AC2" = -0.05618208824592539

    "X3"  = 0.05195879078268548
    "Y3"  = 2.038950300126406
    "AB3" = 0.05216681171301536
    "AC3" = -0.05216681171301536

    "X4"  = 0.05528975755244944
    "Y4"  = 0.4500383411575897
    "AB4" = 0.05502598474440557
    "AC4" = -0.05502598474440557

    "X5"  = 0.05591268251082138
    "Y5"  = 0.1279438197313311
    "AB5" = 0.05618208824592539
    "AC5" = -0.05618208824592539

    "X6"  = 0.06181090156358801
    "Y6"  = 0.1355886256231732
    "AB6" = 0.06169018377745291
    "AC6" = -0.06169018377745291

    "X7"  = 0.06533058382470003
    "Y7"  = 0.6013360828419666
    "AB7" = 0.06313736555146372
    "AC7" = -0.06313736555146372
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
